$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---- Row 11 ----
$ws.Range("A11").Value = "BIOL5081"
$ws.Range("B11").Value = 2
$ws.Range("C11").Value = 43719
$ws.Range("C11").NumberFormat = $ws.Range("C10").NumberFormat
$ws.Range("D11").Value = "the grime reaper"
$ws.Range("D11").NumberFormat = $ws.Range("D10").NumberFormat
$ws.Range("E11").Value = 3
$ws.Range("F11").Value = "won"
$ws.Range("G11").Value = 93
$ws.Range("H11").Value = 4

# ---- Row 12 ----
$ws.Range("A12").Value = "BIOL5081"
$ws.Range("B12").Value = 2
$ws.Range("C12").Value = 43719
$ws.Range("C12").NumberFormat = $ws.Range("C10").NumberFormat
$ws.Range("D12").Value = "the grime reaper"
$ws.Range("D12").NumberFormat = $ws.Range("D10").NumberFormat
$ws.Range("E12").Value = 4
$ws.Range("F12").Value = "won"
$ws.Range("G12").Value = 59
$ws.Range("H12").Value = 5

# ---- Row 13 ----
$ws.Range("A13").Value = "BIOL5081"
$ws.Range("B13").Value = 2
$ws.Range("C13").Value = 43719
$ws.Range("C13").NumberFormat = $ws.Range("C10").NumberFormat
$ws.Range("D13").Value = "the grime reaper"
$ws.Range("D13").NumberFormat = $ws.Range("D10").NumberFormat
$ws.Range("E13").Value = 4
$ws.Range("F13").Value = "won"
$ws.Range("G13").Value = 2
$ws.Range("H13").Value = 3

# Row 13's D:E cells get a font color/name change (matches diff's new font + cellXfs entries)
$ws.Range("D13:E13").Font.Name = "Arial"
$ws.Range("D13:E13").Font.Color = 0

# ---- Selection state to match final workbook view ----
$ws.Range("I13").Select()
